$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the totals row (row 5), shifting the totals row
# and the footer row down by one.
$ws.Rows("5:5").Insert()

# Copy the formatting of the existing product row (row 4) onto the newly
# inserted row 5 so the new line item matches the sheet's styling exactly.
$ws.Range("A4:N4").Copy()
$ws.Range("A5:N5").PasteSpecial(-4122)

# Re-create the merges for the new row (Insert/PasteSpecial above do not
# recreate merged ranges for the brand-new row).
$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()

# Fill in the new product line values.
$ws.Range("A5").Value2 = 2
$ws.Range("B5").Value2 = "مرطب شفاه لونا جوز هند ابيض"
$ws.Range("H5").Value2 = "3:0"
$ws.Range("L5").Value2 = 20
$ws.Range("N5").Value2 = 1

# Update the running total in what is now row 6 (18 + 20 = 38).
$ws.Range("K6").Value2 = 38
